# Add the "Sponsor Covenant" event to the "Confirmation Events" sheet and the
# matching "candidate_events.6.*" columns to the "Candidates with events" sheet.
# (commit message: "17,18: Add Sponsor Covenant & Sponsor Eligibility forms")

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Confirmation Events")
$ws2 = $wb.Worksheets.Item("Candidates with events")

# ---------------------------------------------------------------------------
# Sheet "Confirmation Events": fill in row 8, which was previously blank,
# with a new "Upload Sponsor Covenant" event (same shape/style as row 7).
# ---------------------------------------------------------------------------

# Copy cell formatting (styles only) from the row above so the new row keeps
# the same look (bordered/filled, text vs. date number formats) as the rest
# of the table.
$ws1.Range("A7").Copy()
$ws1.Range("A8").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("B7:C7").Copy()
$ws1.Range("B8:C8").PasteSpecial(-4122)
$ws1.Range("D7").Copy()
$ws1.Range("D8").PasteSpecial(-4122)

$ws1.Range("A8").Value2 = "Upload Sponsor Covenant"
$ws1.Range("B8").Value2 = 42674
$ws1.Range("C8").Value2 = 42658
$ws1.Range("D8").Value2 = "<p><em><strong>Upload Sponsor Covenant</strong></em></p>"

# ---------------------------------------------------------------------------
# Sheet "Candidates with events": add two new columns (Z, AA) holding the
# "candidate_events.6.completed_date" / "candidate_events.6.verified" pair,
# mirroring the existing "candidate_events.5" columns (X, Y).
# ---------------------------------------------------------------------------

# Widen the new columns the same way as the other event columns.
$ws2.Range("Z1").ColumnWidth = $ws2.Range("Y1").ColumnWidth
$ws2.Range("AA1").ColumnWidth = $ws2.Range("Y1").ColumnWidth

# Copy the formatting of the last event pair (X:Y) across the 10 rows onto
# the new columns (Z:AA) so borders / fills / number formats all match.
$ws2.Range("X1:Y10").Copy()
$ws2.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$ws2.Range("Z1").Value2 = "candidate_events.6.completed_date"
$ws2.Range("AA1").Value2 = "candidate_events.6.verified"

# Data rows (2-4 hold real candidates, 5-10 stay blank like the other columns)
$ws2.Range("AA2").Value2 = $false

$ws2.Range("AA3").Value2 = $false

$ws2.Range("Z4").Value2 = 42736
$ws2.Range("AA4").Value2 = $false
